$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: member_id (B6) was stored as a text "1" - convert it to a real number 1
$ws.Range("B6").Value = 1

# Row 7 (new): Amna / Super package payment
$ws.Range("A7").Value = ""
$ws.Range("B7").Value = 4
$ws.Range("C7").Value = ""
$ws.Range("D7").Value = 10000
$ws.Range("E7").Value = ""
$ws.Range("F7").Value = ""
$ws.Range("G7").Value = "Paid"
# keep the date as literal text (not an auto-converted date serial)
$ws.Range("H7").NumberFormat = "@"
$ws.Range("H7").Value = "2025-03-13"
$ws.Range("I7").Value = "Amna"
$ws.Range("J7").Value = "Super"

# Row 8 (new): nabeel / Economy package payment - member_id stays textual "3"
$ws.Range("A8").Value = ""
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "3"
$ws.Range("C8").Value = ""
$ws.Range("D8").Value = 1500
$ws.Range("E8").Value = ""
$ws.Range("F8").Value = ""
$ws.Range("G8").Value = "Paid"
# keep the date as literal text (not an auto-converted date serial)
$ws.Range("H8").NumberFormat = "@"
$ws.Range("H8").Value = "2025-03-13"
$ws.Range("I8").Value = "nabeel"
$ws.Range("J8").Value = "Economy"
